$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.678.74"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.671.92"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.59"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.46"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +4.23%  "
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.399"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.21"
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000196"
$ws.Range("E14").Value = "  -3.91%  "
$ws.Range("D15").Value = "3.151.80"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "65.528.69"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "2.674.05"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.87"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.54"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.71"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000112"
$ws.Range("E24").Value = "  +4.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.63"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("E28").Value = "  -5.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.04"
$ws.Range("E29").Value = "  -5.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "528.64"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.77"
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.49"
$ws.Range("E35").Value = "  -3.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.423"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "157.74"
$ws.Range("E39").Value = "  -3.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.95"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "163.33"
$ws.Range("E42").Value = "  -5.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.13"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0609"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.73"
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0258"
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0264"
$ws.Range("E48").Value = "  +15.00%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.638"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.11"
$ws.Range("E50").Value = "  -4.61%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0995"
$ws.Range("E51").Value = "  +0.22%  "
